$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reverse the order of the "Periodo Mora" values in E16:E21
$ws.Range("E16").Value = "2212"
$ws.Range("E17").Value = "2211"
$ws.Range("E18").Value = "2210"
$ws.Range("E19").Value = "2209"
$ws.Range("E20").Value = "2207"
$ws.Range("E21").Value = "2206"

# Swap the "Valor Mora" amounts between row 16 and row 21
$ws.Range("F16").Value = 35467
$ws.Range("F21").Value = 56000
